# Weekly update: insert a new price record as row 38 ("Fruta / hortaliza, semanal").
# This pushes the existing rows 38-101 down to 39-102 (matching the target
# OOXML diff, which shows every existing record from row 38 onward shifted
# down by exactly one row, plus a brand new row 38 holding the freshly
# reported record).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 38; everything below shifts down one row
# (row 101 -> 102), exactly mirroring the diff's dimension change
# (A1:R101 -> A1:R102).
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new weekly record.
$ws.Range("A38").Value = 2
$ws.Range("B38").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C38").Value = 'Coquimbo'
$ws.Range("D38").Value = 44483
$ws.Range("E38").Value = 4
$ws.Range("F38").Value = 100112031
$ws.Range("G38").Value = 'Poroto verde'
$ws.Range("H38").Value = 'Magnum'
$ws.Range("I38").Value = 'Primera'
$ws.Range("J38").Value = 700
$ws.Range("K38").Value = 43000
$ws.Range("L38").Value = 45000
$ws.Range("M38").Value = 44000
$ws.Range("N38").Value = '$/caja 25 kilos'
$ws.Range("O38").Value = 'Provincia de Limarí'
$ws.Range("P38").Value = 1760
$ws.Range("Q38").Value = 25
$ws.Range("R38").Value = 'Hortaliza'
